$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Post Test" column (C) values for rows 2-23
$values = @(2, 3, 2, 0, 0, 1, 0, 0, 0, 0, 1, 0, 1, 0, 1, 1, 0, 0, 0, 1, 1, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Update the selected cell on the sheet to C24 (one row below the data, column C)
$ws.Range("C24").Select()
